$wb = $excel.ActiveWorkbook

# --- 1. Adjust the existing "fGroups" sheet's selection before it loses focus ---
$ws1 = $wb.Worksheets.Item("fGroups")
$ws1.Activate()
$ws1.Range("B1:G1").Select()

# --- 2. Add the new "mslists" worksheet after fGroups ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "mslists"

# --- 3. Header row (row 1, columns B:G) ---
$ws2.Range("B1").Value = "as-is"
$ws2.Range("C1").Value = "almost as-is"
$ws2.Range("D1").Value = "implement"
$ws2.Range("E1").Value = "not supported"
$ws2.Range("F1").Value = "ionize"
$ws2.Range("G1").Value = "done"

# --- 4. Method rows (A2:G15) ---
$ws2.Range("A2").Value = "`$"
$ws2.Range("B2").Value = "X"

$ws2.Range("A3").Value = "["
$ws2.Range("C3").Value = "X"

$ws2.Range("A4").Value = "[["
$ws2.Range("C4").Value = "X"

$ws2.Range("A5").Value = "analyses"
$ws2.Range("B5").Value = "X"

$ws2.Range("A6").Value = "as.data.table"
$ws2.Range("C6").Value = "X"

$ws2.Range("A7").Value = "averagedPeakLists"
$ws2.Range("C7").Value = "X"

$ws2.Range("A8").Value = "compoundViewer"
$ws2.Range("E8").Value = "X"

$ws2.Range("A9").Value = "filter"
$ws2.Range("C9").Value = "X"

$ws2.Range("A10").Value = "groupNames"
$ws2.Range("B10").Value = "X"

$ws2.Range("A11").Value = "initialize"
$ws2.Range("C11").Value = "X"

$ws2.Range("A12").Value = "length"
$ws2.Range("B12").Value = "X"

$ws2.Range("A13").Value = "peakLists"
$ws2.Range("C13").Value = "X"

$ws2.Range("A14").Value = "plotSpec"
$ws2.Range("B14").Value = "X?"

$ws2.Range("A15").Value = "show"
$ws2.Range("C15").Value = "X"

# --- 5. Styling: column A labels use a monospace "Fira Code" font, vertically centered ---
$rngLabels = $ws2.Range("A2:A15")
$rngLabels.Font.Name = "Fira Code"
$rngLabels.Font.Size = 10
$rngLabels.Font.Color = 0
$rngLabels.VerticalAlignment = -4108

# Last row additionally gets a solid white fill
$ws2.Range("A15").Interior.Color = 16777215

# --- 6. Column widths (best effort, engine quantizes to 1/6-character steps) ---
$ws2.Columns.Item(1).ColumnWidth = 19.833333333333332
$ws2.Columns.Item(3).ColumnWidth = 10.666666666666666
$ws2.Columns.Item(4).ColumnWidth = 10.0
$ws2.Columns.Item(5).ColumnWidth = 12.833333333333334

# --- 7. Page setup ---
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# --- 8. Selection + activation on the new sheet ---
$ws2.Range("D15").Select()
$ws2.Activate()

Write-Output "applied mslists sheet edits"
